$wb = $excel.ActiveWorkbook

# The edit happens on the "credentials" sheet
$ws = $wb.Worksheets.Item("credentials")
$ws.Activate()

# Update the username (B2) and email (C2) values for the first data row
$ws.Range("B2").Value = "u11000"
$ws.Range("C2").Value = "u11000@correo.com"

# Match the author's final cursor/selection position
$ws.Range("C2").Select()
